$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 13538
$ws.Range("I6").Value = 20057
$ws.Range("K6").Value = 60171
$ws.Range("M6").Value = -60059
$ws.Range("H9").Value = 14285857
$ws.Range("I9").Value = 25000140
$ws.Range("J9").Value = 146.66667
$ws.Range("K9").Value = 25000140
$ws.Range("L9").Value = 146.66667
$ws.Range("M9").Value = -24999971
$ws.Range("N9").Value = -484.66667
$ws.Range("H11").Value = 71.86667
$ws.Range("I11").Value = 71.86667
$ws.Range("K11").Value = 71.86667
$ws.Range("M11").Value = 68.13333
$ws.Range("H18").Value = 77368.69500000001
$ws.Range("I18").Value = 91363.73
$ws.Range("J18").Value = 396
$ws.Range("K18").Value = 91363.73
$ws.Range("L18").Value = 396
$ws.Range("M18").Value = -91079.73
$ws.Range("N18").Value = -964
$ws.Range("H32").Value = 2333.2222
$ws.Range("J32").Value = 1928.4286
$ws.Range("L32").Value = 1928.4286
$ws.Range("N32").Value = -2580.4286
$ws.Range("H40").Value = 1673.1562
$ws.Range("I40").Value = 1548.1852
$ws.Range("J40").Value = 2348
$ws.Range("K40").Value = 1548.1852
$ws.Range("L40").Value = 2348
$ws.Range("M40").Value = -1373.1852
$ws.Range("N40").Value = -2698
$ws.Range("H76").Value = 3433.5454
$ws.Range("I76").Value = 2649.5
$ws.Range("J76").Value = 3881.5715
$ws.Range("K76").Value = 2649.5
$ws.Range("L76").Value = 3881.5715
$ws.Range("M76").Value = -2334.5
$ws.Range("N76").Value = -4511.5715
$ws.Range("H79").Value = 3433.5454
$ws.Range("I79").Value = 2649.5
$ws.Range("J79").Value = 3881.5715
$ws.Range("K79").Value = 2649.5
$ws.Range("L79").Value = 3881.5715
$ws.Range("M79").Value = -1557.5
$ws.Range("N79").Value = -6065.5715
$ws.Range("H125").Value = 2752
$ws.Range("I125").Value = 861
$ws.Range("J125").Value = 3292.2856
$ws.Range("K125").Value = 7749
$ws.Range("L125").Value = 29630.5704
$ws.Range("M125").Value = -5289
$ws.Range("N125").Value = -34550.5704
$ws.Range("H129").Value = 948.2857
$ws.Range("I129").Value = 272.33334
$ws.Range("J129").Value = 996
$ws.Range("K129").Value = 817.0000200000001
$ws.Range("L129").Value = 2988
$ws.Range("M129").Value = 4182.99998
$ws.Range("N129").Value = -12988
$ws.Range("H137").Value = 74870.36
$ws.Range("I137").Value = 3596.4
$ws.Range("K137").Value = 10789.2
$ws.Range("M137").Value = -8239.200000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13163369
$ws.Range("I32").Value = 18520470
$ws.Range("J32").Value = 14123.272
$ws.Range("K32").Value = 18520470
$ws.Range("L32").Value = 14123.272
$ws.Range("M32").Value = -18520183
$ws.Range("N32").Value = -14697.272
$ws.Range("H61").Value = 10352.363
$ws.Range("I61").Value = 12357.75
$ws.Range("J61").Value = 5004.6665
$ws.Range("K61").Value = 12357.75
$ws.Range("L61").Value = 5004.6665
$ws.Range("M61").Value = -12145.75
$ws.Range("N61").Value = -5428.6665
$ws.Range("H74").Value = 4437.9355
$ws.Range("I74").Value = 4930.6
$ws.Range("J74").Value = 2385.1667
$ws.Range("K74").Value = 4930.6
$ws.Range("L74").Value = 2385.1667
$ws.Range("M74").Value = -4056.6
$ws.Range("N74").Value = -4133.1667
$ws.Range("H77").Value = 4437.9355
$ws.Range("I77").Value = 4930.6
$ws.Range("J77").Value = 2385.1667
$ws.Range("K77").Value = 24653
$ws.Range("L77").Value = 11925.8335
$ws.Range("M77").Value = -20285
$ws.Range("N77").Value = -20661.8335
$ws.Range("H124").Value = 20785.9
$ws.Range("J124").Value = 20785.9
$ws.Range("L124").Value = 20785.9
$ws.Range("N124").Value = -30605.9
$ws.Range("H136").Value = 10352.363
$ws.Range("I136").Value = 12357.75
$ws.Range("J136").Value = 5004.6665
$ws.Range("K136").Value = 37073.25
$ws.Range("L136").Value = 15013.9995
$ws.Range("M136").Value = -34523.25
$ws.Range("N136").Value = -20113.9995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 5751.1665
$ws.Range("I11").Value = 375.5
$ws.Range("J11").Value = 16502.5
$ws.Range("K11").Value = 375.5
$ws.Range("L11").Value = 16502.5
$ws.Range("M11").Value = -235.5
$ws.Range("N11").Value = -16782.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 340.27274
$ws.Range("I22").Value = 333.22223
$ws.Range("J22").Value = 372
$ws.Range("K22").Value = 333.22223
$ws.Range("L22").Value = 372
$ws.Range("M22").Value = 16.77776999999998
$ws.Range("N22").Value = -1072
$ws.Range("H134").Value = 8075.75
$ws.Range("I134").Value = 2800
$ws.Range("J134").Value = 8829.429
$ws.Range("K134").Value = 8400
$ws.Range("L134").Value = 26488.287
$ws.Range("M134").Value = -5865
$ws.Range("N134").Value = -31558.287

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 40000028
$ws.Range("I33").Value = 18.333334
$ws.Range("J33").Value = 100000050
$ws.Range("K33").Value = 110.000004
$ws.Range("L33").Value = 600000300
$ws.Range("M33").Value = 172.999996
$ws.Range("N33").Value = -600000866
$ws.Range("H68").Value = 1050.1702
$ws.Range("I68").Value = 969.4167
$ws.Range("J68").Value = 1134.4348
$ws.Range("K68").Value = 2908.2501
$ws.Range("L68").Value = 3403.3044
$ws.Range("M68").Value = -2097.2501
$ws.Range("N68").Value = -5025.3044
$ws.Range("H71").Value = 1050.1702
$ws.Range("I71").Value = 969.4167
$ws.Range("J71").Value = 1134.4348
$ws.Range("K71").Value = 8724.7503
$ws.Range("L71").Value = 10209.9132
$ws.Range("M71").Value = -4668.7503
$ws.Range("N71").Value = -18321.9132

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 943.55554
$ws.Range("I16").Value = 943.55554
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 943.55554
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -773.55554
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 501210
$ws.Range("I22").Value = 667828
$ws.Range("K22").Value = 667828
$ws.Range("M22").Value = -667533
$ws.Range("H27").Value = 501210
$ws.Range("I27").Value = 667828
$ws.Range("K27").Value = 667828
$ws.Range("M27").Value = -667721
$ws.Range("H46").Value = 250001680
$ws.Range("I46").Value = 500000740
$ws.Range("J46").Value = 2600
$ws.Range("K46").Value = 500000740
$ws.Range("L46").Value = 2600
$ws.Range("M46").Value = -500000552
$ws.Range("N46").Value = -2976
$ws.Range("H55").Value = 1443658.6
$ws.Range("I55").Value = 2525652.5
$ws.Range("K55").Value = 2525652.5
$ws.Range("M55").Value = -2525479.5
$ws.Range("H127").Value = 31928.75
$ws.Range("J127").Value = 31928.75
$ws.Range("L127").Value = 31928.75
$ws.Range("N127").Value = -41848.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 1670816.6
$ws.Range("I7").Value = 3333700
$ws.Range("J7").Value = 7933.3335
$ws.Range("K7").Value = 3333700
$ws.Range("L7").Value = 7933.3335
$ws.Range("M7").Value = -3333587
$ws.Range("N7").Value = -8159.3335
$ws.Range("H132").Value = 1660.9714
$ws.Range("I132").Value = 1381.24
$ws.Range("K132").Value = 4143.72
$ws.Range("M132").Value = -1613.72
$ws.Range("H136").Value = 37690.43
$ws.Range("I136").Value = 63613.25
$ws.Range("K136").Value = 190839.75
$ws.Range("M136").Value = -188289.75

Write-Output "edits applied"